# Adding the changes we made on may 9th
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to insert right after the header (row 1), shifting existing data down.
$newTop = @(
    @(-0.1574550628662107, 0.1197370529174804, -0.3381910085678101),
    @(-0.3228707909584047, 0.1800707578659058, -0.2396511320024727),
    @(-0.3687732219696044, 0.22406005859375, -0.1331830434501169),
    @(-0.3819821834564209, 0.1978947639465331, -0.2025355234742165),
    @(-0.3319910526275633, 0.2015152931213379, -0.2596112400293351),
    @(-0.3998157501220704, 0.1588554024696349, -0.4224489629268647),
    @(-0.4262521266937255, -0.03637795448303255, -0.570676848292351),
    @(-0.5638278126716616, -0.07290470600128157, -1.039905160665513)
)

# Insert 8 new rows before existing row 2 (pushing old data down).
$insertRange = $ws.Range("A2:C9")
$insertRange.Insert()
$insertRange.ClearFormats()

for ($i = 0; $i -lt $newTop.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newTop[$i][0]
    $ws.Cells.Item($r, 2).Value = $newTop[$i][1]
    $ws.Cells.Item($r, 3).Value = $newTop[$i][2]
}

# New rows to append at the end (rows 30 and 31).
$newBottom = @(
    @(-0.003773307800292898, 0.2548617899417878, -0.1365976125001907),
    @(0.01674773693084754, 0.2881045341491697, -0.1470022119581697)
)

for ($i = 0; $i -lt $newBottom.Count; $i++) {
    $r = 30 + $i
    $ws.Cells.Item($r, 1).Value = $newBottom[$i][0]
    $ws.Cells.Item($r, 2).Value = $newBottom[$i][1]
    $ws.Cells.Item($r, 3).Value = $newBottom[$i][2]
}
